$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 7 with the new data for the "generateReceipt" task
$ws.Range("B7").Value = "10 mins"
$ws.Range("C7").Value = "7 mins 18 seconds"
$ws.Range("D7").Value = "have only created pseudocode and added the method for the sub-tasks all the while setting the receipt format"
$ws.Range("E7").Value = "-"

# Row 7 needs to be taller to match the wrapped text rows (like rows 2-4)
$ws.Rows.Item(7).RowHeight = 45

# Update the selection to match the new active cell / range
$ws.Range("B8:E8").Select()
